$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to include the extra guidance text shown on the
# live site's data-entry form.
$ws.Range("D1").Value = "Is Virtual Machine (TRUE / FALSE)"
$ws.Range("E1").Value = "Environment (Prod / Dev / QA)"
$ws.Range("V1").Value = "Start Date (MM/DD/YY)"
$ws.Range("W1").Value = "Next Hardware Support Date (MM/DD/YY)"
$ws.Range("X1").Value = "Base Warranty (MM/DD/YY)"

# Widen the columns that now hold longer header text so it matches the
# site's display formatting. (Values are pre-compensated for the COM
# layer's column-width quantization so the saved OOXML width lands on
# the intended target.)
$ws.Columns.Item(4).ColumnWidth = 37.166666666666664
$ws.Columns.Item(5).ColumnWidth = 26.833333333333332
$ws.Columns.Item(23).ColumnWidth = 36.5
$ws.Columns.Item(24).ColumnWidth = 26.666666666666668

# Move the active selection to match the new view state.
$ws.Range("G4").Select() | Out-Null
